$p = $ppt.ActivePresentation

# The deck's fixed "datetime1" / "datetimeFigureOut" footer fields (on the
# slide master, all 11 slide layouts, and the notes master) are being bumped
# from 12/11/2024 to 2/14/2025. Each lives in a "Date Placeholder" shape.
$newDate = "2/14/2025"

$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$layoutDateShapeIndex = @{
    1 = 3   # Title Slide
    2 = 3   # Title and Content
    3 = 3   # Section Header
    4 = 4   # Two Content
    5 = 6   # Comparison
    6 = 2   # Title Only
    7 = 1   # Blank
    8 = 4   # Content with Caption
    9 = 4   # Picture with Caption
    10 = 3  # Title and Vertical Text
    11 = 3  # Vertical Title and Text
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $shapeIdx = $layoutDateShapeIndex[$i]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = $newDate
}

$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Slide 1's title placeholder ("Title 1") moved from (1524000, 1840342) EMU to
# (1689100, 239886) EMU; size is unchanged. Shape.Left/Top are in points
# (1 pt = 12700 EMU).
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(1)
$title.Left = 1689100 / 12700.0
$title.Top = 239886 / 12700.0
